# Auto-generated Excel COM-interop edit script
# Applies the weekly CompStat report refresh: updates volume/date header text,
# and refreshes all crime-stat figures in rows 14-33 (values, and some cells
# switching between numeric/percentage and "no data" (N/A) text placeholders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: volume number and report week dates (rich-text runs) ---
$ws.Range("A8").Characters(21, 2).Text = "26"
$ws.Range("C9").Characters(27, 9).Text = "6/23/2025"
$ws.Range("C9").Characters(47, 9).Text = "6/29/2025"

# --- Plain value updates (style/number-format unchanged) ---
$ws.Range("N14").Value = -90.909090909090
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("L15").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 5
$ws.Range("I16").Value = 41
$ws.Range("J16").Value = 24
$ws.Range("K16").Value = 70.833333333333
$ws.Range("L16").Value = -14.583333333333
$ws.Range("M16").Value = -74.050632911392
$ws.Range("N16").Value = -91.700404858299
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -57.142857142857
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -52.941176470588
$ws.Range("I17").Value = 103
$ws.Range("J17").Value = 137
$ws.Range("K17").Value = -24.817518248175
$ws.Range("L17").Value = -11.965811965812
$ws.Range("M17").Value = -31.788079470198
$ws.Range("N17").Value = -44.021739130434
$ws.Range("C18").Value = 4
$ws.Range("F18").Value = 9
$ws.Range("H18").Value = 80
$ws.Range("I18").Value = 68
$ws.Range("J18").Value = 54
$ws.Range("K18").Value = 25.925925925925
$ws.Range("L18").Value = -2.857142857142
$ws.Range("M18").Value = -61.581920903954
$ws.Range("N18").Value = -91.468005018820
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -40
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 178
$ws.Range("J19").Value = 185
$ws.Range("K19").Value = -3.783783783783
$ws.Range("L19").Value = -9.183673469387
$ws.Range("M19").Value = -19.090909090909
$ws.Range("N19").Value = -39.864864864864
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -33.333333333333
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 8.333333333333
$ws.Range("I20").Value = 91
$ws.Range("J20").Value = 117
$ws.Range("K20").Value = -22.222222222222
$ws.Range("L20").Value = 15.189873417721
$ws.Range("M20").Value = -50.543478260869
$ws.Range("N20").Value = -94.547633313361
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -31.818181818181
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = -4.838709677419
$ws.Range("I21").Value = 489
$ws.Range("J21").Value = 527
$ws.Range("K21").Value = -7.210626185958
$ws.Range("L21").Value = -5.598455598455
$ws.Range("M21").Value = -46.263736263736
$ws.Range("N21").Value = -85.924006908462
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 142.857142857143
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 30.434782608695
$ws.Range("I24").Value = 280
$ws.Range("J24").Value = 347
$ws.Range("K24").Value = -19.308357348703
$ws.Range("L24").Value = -28.571428571428
$ws.Range("M24").Value = -33.962264150943
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 14
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -17.647058823529
$ws.Range("I25").Value = 71
$ws.Range("J25").Value = 85
$ws.Range("K25").Value = -16.470588235294
$ws.Range("L25").Value = -1.388888888888
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = -56.25
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 53
$ws.Range("H26").Value = -41.509433962264
$ws.Range("I26").Value = 217
$ws.Range("J26").Value = 205
$ws.Range("K26").Value = 5.853658536585
$ws.Range("L26").Value = 26.162790697674
$ws.Range("M26").Value = -32.398753894081
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 0
$ws.Range("L27").Value = -36.363636363636
$ws.Range("J29").Value = 3
$ws.Range("J30").Value = 3
$ws.Range("I33").Value = 2

# --- Cells switching from numeric/percentage to "N/A" text placeholders ---
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$ws.Range("E15").Value = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$ws.Range("C16").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C16").PasteSpecial($xlPasteFormats)
$ws.Range("D27").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$ws.Range("E27").Value = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial($xlPasteFormats)
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial($xlPasteFormats)
$ws.Range("E28").Value = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial($xlPasteFormats)
$ws.Range("F31").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F31").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# --- Cells switching from "N/A" text placeholders to numeric/percentage values ---
$ws.Range("D16").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D16").PasteSpecial($xlPasteFormats)
$ws.Range("E16").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E16").PasteSpecial($xlPasteFormats)
$ws.Range("D18").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial($xlPasteFormats)
$ws.Range("E18").Value = 300
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("D29").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial($xlPasteFormats)
$ws.Range("E29").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial($xlPasteFormats)
$ws.Range("G29").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("G29").PasteSpecial($xlPasteFormats)
$ws.Range("H29").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("H29").PasteSpecial($xlPasteFormats)
$ws.Range("D30").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial($xlPasteFormats)
$ws.Range("E30").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial($xlPasteFormats)
$ws.Range("G30").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("G30").PasteSpecial($xlPasteFormats)
$ws.Range("H30").Value = -100
$ws.Range("L14").Copy() | Out-Null
$ws.Range("H30").PasteSpecial($xlPasteFormats)
$ws.Range("C33").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C33").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
